$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "related_words" column (D) for several quality attributes to
# tightened/refreshed keyword lists.
$ws.Range("D2").Value = "{stability,fault tolerance,robustness,trustworthiness}"
$ws.Range("D4").Value = "{authentication,authorization,auth,login}"
$ws.Range("D12").Value = "{type of license,license compliance requirements,licensing,proprietary,copyright,copyleft,usage limits}"
$ws.Range("D15").Value = "{modularity,decoupling,cohesion,robustness,observability,controllability,verifiability,diagnosability,unit test,test driven,test,examin,evaluate}"
$ws.Range("D19").Value = "{compatible,interoperability,coexistence,integration,conformance,alignment}"

# Leave the cursor where the author last left it.
$ws.Range("D5").Select()
